# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The workers table (Hoja1!B15:J22) is reshuffled: the rows for
# "CRISTIAN DE JESUS TAPIA BERDUGO" (doc 1143405636) and
# "ROLANDO JAVIER RAMOS BURGOS" (doc 73187252) are interleaved by period,
# and CRISTIAN's overdue amount (Valor Mora) is updated from 1300000 to
# 2000000 for his three periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: CRISTIAN - period 1811 - 2,000,000
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143405636"
$ws.Range("D16").Value = "CRISTIAN DE JESUS TAPIA BERDUGO"
$ws.Range("E16").Value = "1811"
$ws.Range("F16").Value = 31249
$ws.Range("G16").Value = 2000000

# Row 17: ROLANDO - period 1812 - 877,803
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73187252"
$ws.Range("D17").Value = "ROLANDO JAVIER RAMOS BURGOS"
$ws.Range("E17").Value = "1812"
$ws.Range("F17").Value = 31249
$ws.Range("G17").Value = 877803

# Row 18: CRISTIAN - period 1812 - 2,000,000
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143405636"
$ws.Range("D18").Value = "CRISTIAN DE JESUS TAPIA BERDUGO"
$ws.Range("E18").Value = "1812"
$ws.Range("F18").Value = 31249
$ws.Range("G18").Value = 2000000

# Row 19: ROLANDO - period 1901 - 877,803
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73187252"
$ws.Range("D19").Value = "ROLANDO JAVIER RAMOS BURGOS"
$ws.Range("E19").Value = "1901"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 877803

# Row 20: CRISTIAN - period 1901 - 2,000,000
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143405636"
$ws.Range("D20").Value = "CRISTIAN DE JESUS TAPIA BERDUGO"
$ws.Range("E20").Value = "1901"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 2000000

# Row 21: ROLANDO - period 1902 - 877,803
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "73187252"
$ws.Range("D21").Value = "ROLANDO JAVIER RAMOS BURGOS"
$ws.Range("E21").Value = "1902"
$ws.Range("F21").Value = 31249
$ws.Range("G21").Value = 877803

# Row 22: ROLANDO - period 1907 - 877,803
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "73187252"
$ws.Range("D22").Value = "ROLANDO JAVIER RAMOS BURGOS"
$ws.Range("E22").Value = "1907"
$ws.Range("F22").Value = 31249
$ws.Range("G22").Value = 877803
